$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("L1").Value = "saldo pokok"
$ws.Range("M1").Value = "saldo wajib"
$ws.Range("N1").Value = "saldo manasuka"

$ws.Range("L1").Select()
